# Fixes #4 and Fixes #5 . Added gerber 1.4.1 for hotplate with corrected
# playout and removed 1.4 and fixed Inductor note in the BOM.
#
# The only real data change in this BOM workbook is row 23 (the L1
# inductor line): the part was switched from the LCSC-sourced
# "CDRH74-220MT" / MetalLions part to an AliExpress-sourced
# "CDRH74-221MT" (220uH, CD74R-221) part, with the Manufacturer cell
# turned into a hyperlink to the AliExpress listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aliexpressUrl = "https://www.aliexpress.com/item/1005005453639618.html?spm=a2g0o.productlist.main.7.fc121f071PKJLW&algo_pvid=f9629d7f-965b-4802-a78f-4d249e8c6e80&algo_exp_id=f9629d7f-965b-4802-a78f-4d249e8c6e80-3&pdp_npi=4%40dis%21BGN%214.20%213.99%21%21%2116.56%2115.73%21%402101c5a417171395910791417e39c0%2112000033143312160%21sea%21BG%213156952160%21&curPageLogUid=4Y1aJAra4aBs&utparam-url=scene%3Asearch%7Cquery_from%3A"

# G23: Manufacturer Part
$ws.Range("G23").Value = "220uH (CD74R-221)"

# I23: Supplier Part
$ws.Range("I23").Value = "AliExpress"

# J23: Supplier
$ws.Range("J23").Value = "N/A"

# N23: Device
$ws.Range("N23").Value = "CDRH74-221MT"

# H23: Manufacturer -> becomes a hyperlink to the AliExpress listing,
# with the link text/display set to the URL itself (also applies the
# built-in Hyperlink style: underlined, theme colour 10).
$ws.Hyperlinks.Add($ws.Range("H23"), $aliexpressUrl, "", "", $aliexpressUrl) | Out-Null

# View state: zoomed to 70% with M37 selected (matches the author's
# session at save time).
$excel.ActiveWindow.Zoom = 70
$ws.Range("M37").Select() | Out-Null
